$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the number formatting/style from column Q (rows 4-14) into the new
# column R so the added cells inherit the same look (year header, data rows,
# bottom border on the last row) as the rest of the table.
$ws.Range("Q4:Q14").Copy()
$ws.Range("R4:R14").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# New "2020" column values, mirroring the existing 2007..2019 columns.
$ws.Range("R4").Value = 2020
$ws.Range("R5").Value = 2.1
$ws.Range("R6").Value = 2.4
$ws.Range("R7").Value = 1.4
$ws.Range("R8").Value = 3.2
$ws.Range("R9").Value = 2.4
$ws.Range("R10").Value = 0.8
$ws.Range("R11").Value = 2.2000000000000002
$ws.Range("R12").Value = 4.5
$ws.Range("R13").Value = 1.4
$ws.Range("R14").Value = 3.2

# Move/extend the active selection the way the author left it.
$ws.Range("R16:R17").Select()
